# "Add third (setup) button"
# Draws the bitmap for the 3rd extra custom LCD character group (columns S:W,
# the "d"/"h" characters) and finishes the "g" character (columns M:Q, bottom
# block), then marks the newly-drawn top block (S3:W5) with its own
# "highlight when 1" conditional-format rule, matching the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gun")

# --- custom_char_d bitmap (rows 3-5, cols S:W) -----------------------------
# row3 = B10110, row4 = B11111, row5 = B10110
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 1
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = 0

$ws.Range("S4:W4").Value = 1

$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 1
$ws.Range("V5").Value = 1
$ws.Range("W5").Value = 0

# --- custom_char_g bitmap (rows 14-16, cols M:Q) ---------------------------
# row14 = B10110, row15 = B11111, row16 = B10110
$ws.Range("M14").Value = 1
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 1
$ws.Range("P14").Value = 1
$ws.Range("Q14").Value = 0

$ws.Range("M15:Q15").Value = 1

$ws.Range("M16").Value = 1
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 1
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 0

# --- custom_char_h bitmap (rows 14-16, cols S:W) ---------------------------
# row14 = B11100, row15 = B11110, row16 = B11100
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 1
$ws.Range("U14").Value = 1
$ws.Range("V14").Value = 0
$ws.Range("W14").Value = 0

$ws.Range("S15").Value = 1
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = 1
$ws.Range("V15").Value = 1
$ws.Range("W15").Value = 0

$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1
$ws.Range("U16").Value = 1
$ws.Range("V16").Value = 0
$ws.Range("W16").Value = 0

# --- conditional formatting: split S3:W10 into S3:W5 (new rule, top
#     priority) + S6:W10 (existing rule, shrunk) so the freshly-drawn block
#     gets its own "highlight when 1" rule like every other glyph block ----
$oldRule = $ws.Range("S3:W10").FormatConditions.Item(1)
$oldRule.ModifyAppliesToRange($ws.Range("S6:W10"))

$newRule = $ws.Range("S3:W5").FormatConditions.Add(1, 3, 1)
$newRule.Font.Color = 393372
$newRule.Interior.Color = 13551615
$newRule.SetFirstPriority()
